# Update loading_percent results for the 380 kV case (Case_1_81 / res_line).
# Columns affected: B, D, E, F, G, J, K, N, O across rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "E", "F", "G", "J", "K", "N", "O")

# Each inner array: row number followed by the new value for each column in $cols, in order.
$data = @(
    @(2, 7.802163441331366, 9.148204694649557, 13.66263906284698, 32.98997795880901, 3.652609188650061, 9.926579855817202, 12.33396464199977, 18.71970345397722, 24.81047487196775),
    @(3, 7.73022618312397, 9.097931103324527, 13.60361568626759, 32.99979720809023, 3.654714673885186, 9.93266694723243, 11.99242327940801, 18.7781689841043, 24.86407482211169),
    @(4, 7.687554255760942, 9.068507454728595, 13.57012779578047, 33.01421444705095, 3.656076423203922, 9.938011878245279, 11.77926643092534, 18.81593448140726, 24.9027067789354),
    @(5, 7.670561217975227, 9.056889759720352, 13.55718391833187, 33.0221974786387, 3.656648746368302, 9.94059447789787, 11.69166443551858, 18.83179501806228, 24.91988482883057),
    @(6, 7.667764016277652, 9.054983424924377, 13.55507732896185, 33.0236503019979, 3.656744832698577, 9.941047753587906, 11.67707747176043, 18.83445711841104, 24.92282380637759),
    @(7, 7.68732345264837, 9.0683492533723, 13.56995037164437, 33.01431357701284, 3.656084071230856, 9.938045069990416, 11.77808781593862, 18.81614647433333, 24.90293264186154),
    @(8, 7.777060499202483, 9.130576362825622, 13.64172306024211, 32.99162218961269, 3.653320878081886, 9.928345270947981, 12.21699165247838, 18.73947551364553, 24.82776667252026),
    @(9, 7.96403102391786, 9.263627959934679, 13.80380808955783, 33.01369295144377, 3.648447003213953, 9.922061719780418, 13.04508089975344, 18.60388575673616, 24.72590556323273),
    @(10, 8.106875603111568, 9.367486032454048, 13.93517261049212, 33.07043352919422, 3.645194714018115, 9.925182903708169, 13.62719721490761, 18.51318931181699, 24.67900602901208),
    @(11, 8.172781065334027, 9.41592014076123, 13.99742923791245, 33.10500690471713, 3.643785751207516, 9.928274639623373, 13.8851296905974, 18.47384994130229, 24.66376661790119),
    @(12, 8.197848132610572, 9.434419557784722, 14.02134796294409, 33.11935325957197, 3.643262297196658, 9.929684863578471, 13.98172659513868, 18.45922787958336, 24.65887410168882),
    @(13, 8.192444957904879, 9.430428529971135, 14.01618161777795, 33.11620784224205, 3.643374584452856, 9.929370513779634, 13.96097191293032, 18.46236479197988, 24.65988870508636),
    @(14, 8.174841272269026, 9.417438995067004, 13.99939025407191, 33.10616209045651, 3.643742484410453, 9.928385866653482, 13.89309882499052, 18.4726414718361, 24.66334649606827),
    @(15, 8.164072166878606, 9.40950279404505, 13.98914931399019, 33.10017191264588, 3.643969146216391, 9.927813894470388, 13.85138193770713, 18.47897200533111, 24.66557892190676),
    @(16, 8.102585177525693, 9.364343608652138, 13.93115293176564, 33.06834992567025, 3.645288207719006, 9.9250144147701, 13.61019467064077, 18.51579875676033, 24.68012477995453),
    @(17, 8.065085461184619, 9.336935559597256, 13.89620287716569, 33.05106860387475, 3.646115434607596, 9.923724591987325, 13.46040799254729, 18.53888151210914, 24.69061079748256),
    @(18, 8.043605027850681, 9.321283809306763, 13.87633666070075, 33.04195397876047, 3.646597874322919, 9.92314015279565, 13.37361224450511, 18.55233877858357, 24.69721572660242),
    @(19, 8.036348004819125, 9.316004115478579, 13.86965132681357, 33.03900979435758, 3.646762362227592, 9.922969338196316, 13.34411717649241, 18.55692624081062, 24.69955050489196),
    @(20, 8.069068385641778, 9.339841625260755, 13.89989904372777, 33.05282286922289, 3.646026687977729, 9.923845608696846, 13.47642019064516, 18.53640562161531, 24.68943515800222),
    @(21, 8.180009095353686, 9.421250135842344, 14.00431309084255, 33.10907878623462, 3.643634149850485, 9.928668591116711, 13.91306467181946, 18.46961550915687, 24.66230700952762),
    @(22, 8.253146472107332, 9.475373062622392, 14.07454812836917, 33.15315189491129, 3.642129274541782, 9.933215848585874, 14.19212731282397, 18.42756619603274, 24.64969723207649),
    @(23, 8.214061454579538, 9.446406834133734, 14.0368852308809, 33.12896294884612, 3.642927092593154, 9.930661587943321, 14.04379056799654, 18.44986246892663, 24.65595834149883),
    @(24, 8.067267458764439, 9.338527463413529, 13.89822729845872, 33.05202720917076, 3.646066788968826, 9.923790407629731, 13.46918319148028, 18.53752439011024, 24.68996486933263),
    @(25, 7.912398388428113, 9.226515639311012, 13.75774730629466, 33.00059626989763, 3.649707564968808, 9.922400269593281, 13.04508089975344, 18.63899414259977, 24.72590556323273)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $col = $cols[$i]
        $val = $row[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}
